# Generate Report for Handoff
# Adds a new localization-status row (for c63cf7c5-3c74-4347-89c7-bb8a5cf6fb9c.md)
# to each of the three report sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/760380b21233d6f346233f88ffb9952079e1d09a/e2e"

$newFile        = "c63cf7c5-3c74-4347-89c7-bb8a5cf6fb9c.md"
$newFileDisplay = "e2e\c63cf7c5-3c74-4347-89c7-bb8a5cf6fb9c.md"
$newFileUrl     = "$repoBase/$newFile"

$hyperlinkColor = 15570276  # BGR long for RGB FF6495ED (matches the workbook's existing HyperLink font colour)

# A leading apostrophe forces text (not bool/number) interpretation for
# literal "True"/"False"/"" values, matching the source report's columns
# which are all typed as shared strings, never booleans.
$emptyText = "'"
$trueText  = "'True"
$falseText = "'False"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$lo1 = $ws1.ListObjects.Item(1)
$lo1.ListRows.Add() | Out-Null

$ws1.Range("A3").Value = $newFile
$ws1.Range("B3").Value = $newFileDisplay
$ws1.Range("C3").Value = ".md"
$ws1.Range("D3").Value = $emptyText
$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-08-17 08:39:20"
$ws1.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws1.Hyperlinks.Add($ws1.Range("B3"), $newFileUrl, "", "", $newFileDisplay) | Out-Null
$ws1.Range("B3").Font.Underline = 2
$ws1.Range("B3").Font.Color = $hyperlinkColor

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$lo2 = $ws2.ListObjects.Item(1)
$lo2.ListRows.Add() | Out-Null

$ws2.Range("A3").Value = $newFile
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "e2e"
$ws2.Range("E3").Value = "ht"
$ws2.Range("F3").Value = $falseText
$ws2.Range("G3").Value = "c63cf7c5-3c74-4347-89c7-bb8a5cf6fb9c.d9189d51265ff6c38dd8b28b9bbc00ea754b81f2.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-17 08:39:15"
$ws2.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("I3").Value = $emptyText
$ws2.Range("J3").Value = $emptyText
$ws2.Range("K3").Value = "0001-01-01 00:00:00"
$ws2.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("L3").Value = $emptyText
$ws2.Range("M3").Value = $trueText
$ws2.Range("N3").Value = $emptyText
$ws2.Range("O3").Value = $falseText
$ws2.Range("P3").Value = $emptyText

$ws2.Hyperlinks.Add($ws2.Range("A3"), $newFileUrl, "", "", $newFile) | Out-Null
$ws2.Range("A3").Font.Underline = 2
$ws2.Range("A3").Font.Color = $hyperlinkColor

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$lo3 = $ws3.ListObjects.Item(1)
$lo3.ListRows.Add() | Out-Null

$ws3.Range("A3").Value = $newFile
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "e2e"
$ws3.Range("E3").Value = "ht"
$ws3.Range("F3").Value = $falseText
$ws3.Range("G3").Value = "c63cf7c5-3c74-4347-89c7-bb8a5cf6fb9c.d9189d51265ff6c38dd8b28b9bbc00ea754b81f2.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-17 08:39:20"
$ws3.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("I3").Value = $emptyText
$ws3.Range("J3").Value = $emptyText
$ws3.Range("K3").Value = "0001-01-01 00:00:00"
$ws3.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("L3").Value = $emptyText
$ws3.Range("M3").Value = $trueText
$ws3.Range("N3").Value = $emptyText
$ws3.Range("O3").Value = $falseText
$ws3.Range("P3").Value = $emptyText

$ws3.Hyperlinks.Add($ws3.Range("A3"), $newFileUrl, "", "", $newFile) | Out-Null
$ws3.Range("A3").Font.Underline = 2
$ws3.Range("A3").Font.Color = $hyperlinkColor

Write-Host "Handoff report row added to Overview, zh-cn, de-de."
